# Reflow the lyric text across several slides so the text wraps by
# character-length instead of whole-word boundaries (per commit
# "using length instead of words").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Song "Con vung tin noi Ngai" (slides 7-12) + "Khi Nghi Ve Ngai"
# divider/refrain (slides 12-14 before the insert).
#
# A brand new slide is needed to hold the text that overflows once
# the wrap points move, so duplicate the existing lyric slide 11
# (same layout/formatting) - this inserts the copy right after it,
# at position 12, and pushes everything from the old slide 12 on
# down by one.
# ---------------------------------------------------------------
$dup = $p.Slides.Item(11).Duplicate()

# Re-flowed verse 2 / verse 3 lyric lines (slides 7-12).
$p.Slides.Item(7).Shapes.Item(1).TextFrame.TextRange.Text = "Nguyện dâng lên Chúa những tin yêu. Niềm tin nhỏ bé rất chênh vênh, rất lung lay Chúa con ơi! Xin Ngài xót"
$p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange.Text = "thương."
$p.Slides.Item(9).Shapes.Item(1).TextFrame.TextRange.Text = "2. Con vững tin nơi Ngài là nguồn ánh sáng chiếu đời con. Con vững tin nơi Ngài là đường để con bước đi. Con"
$p.Slides.Item(10).Shapes.Item(1).TextFrame.TextRange.Text = "vững tin nơi Ngài dù bóng tôi khuất lối đi, dù mây đen giăng mịt mù dựa vào Chúa dẫn đường bước đi."
$p.Slides.Item(11).Shapes.Item(1).TextFrame.TextRange.Text = "3. Con vững tin nơi Ngài là dòng suối tưới mát hồn con. Con vững tin nơi Ngài tình Ngài đại dương xót thương."
$p.Slides.Item(12).Shapes.Item(1).TextFrame.TextRange.Text = "Con vững tin nơi Ngài tội lỗi dẫu có chất chồng dù con bao phen lạc đường một lòng tín thác tình Chúa thương."

# The old slide 12 (the "Nhap Le / Khi Nghi Ve Ngai" divider) is now
# at position 13, unchanged. The old slide 13 ("DK: Khi nghi ve
# Ngai...") is now at position 14, and the old slide 14 ("Ngai.") is
# now at position 15 - its single word is folded back into slide 14
# and the now-empty slide is removed.
$p.Slides.Item(14).Shapes.Item(1).TextFrame.TextRange.Text = "ĐK: Khi nghĩ về Ngài, con tự nhủ: Hãy đi tìm thánh nhan. Đừng ẩn xa con ôi lạy Chúa, con đi tìm thánh nhan Ngài."
$p.Slides.Item(15).Delete()

# ---------------------------------------------------------------
# Song "Chua Bien Hinh" refrain (slides 18-19): move "vang." from
# the end of slide 18 to the start of slide 19.
# ---------------------------------------------------------------
$p.Slides.Item(18).Shapes.Item(1).TextFrame.TextRange.Text = "ĐK: Hãy biến đổi xa lợi danh thế trần. Sống với Chúa tâm hồn sẽ thanh nhàn. Chúa giúp sức cho niềm tin vững"
$p.Slides.Item(19).Shapes.Item(1).TextFrame.TextRange.Text = "vàng. Hãy giũ hết bao tội lỗi đã mang."

# ---------------------------------------------------------------
# Song "Dang Tron Cuoc Doi" (slides 21-22): move "huong" from the
# end of slide 21 to the start of slide 22. Slide 21's text contains
# an embedded line break inside a single paragraph/run, so trim the
# trailing word via Characters() instead of reassigning the whole
# TextRange.Text (which would split it into two paragraphs).
# ---------------------------------------------------------------
$tr21 = $p.Slides.Item(21).Shapes.Item(1).TextFrame.TextRange
$tr21.Characters($tr21.Length - 5, 6).Text = ""
$p.Slides.Item(22).Shapes.Item(1).TextFrame.TextRange.Text = "hương hoa vườn trái phương xa. Hiệp lòng dâng Cha cùng với câu ca. Và niềm mơ ước đời sống an hòa."

# ---------------------------------------------------------------
# Song "Ta On Chua Xuan" refrain (slides 27-28): move "ta on" from
# the end of slide 27 to the start of slide 28.
# ---------------------------------------------------------------
$p.Slides.Item(27).Shapes.Item(1).TextFrame.TextRange.Text = "ĐK: Đây Chúa đem mùa xuân, Chúa đem an hòa, trải rộng trần gian. Muôn trái tim hiệp thông, thiết tha dâng lời,"
$p.Slides.Item(28).Shapes.Item(1).TextFrame.TextRange.Text = "tạ ơn Chúa Xuân."
